$d = $word.ActiveDocument

$pairs = @(
    @("779×6=", "665×6="),
    @("577×2=", "205×4="),
    @("310×5=", "182×8="),
    @("826×8=", "244×8="),
    @("266×5=", "861×2="),
    @("914×7=", "400×4="),
    @("301×5=", "883×3="),
    @("530×6=", "520×6="),
    @("133×6=", "412×9="),
    @("574×5=", "227×3="),
    @("449×8=", "380×5="),
    @("878×9=", "432×2="),
    @("157×9=", "750×9="),
    @("664×4=", "854×7="),
    @("453×8=", "517×7="),
    @("242×6=", "517×5="),
    @("828×7=", "603×8="),
    @("583×8=", "179×3="),
    @("364×5=", "361×7="),
    @("843×9=", "789×2="),
    @("650×8=", "250×3="),
    @("567×8=", "536×7="),
    @("655×8=", "948×6="),
    @("507×7=", "603×2="),
    @("650×6=", "792×3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
